$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed LR-pair data (updated TPM values), including a new
# "Inflammatory-Mac" sending-cluster block inserted after "FAPs" and
# before "MuSCs", matching the regenerated NATMI output.
$data = @(
    @("ECs", "Tgfb2", "Acvr1", "ECs", 3, 1, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3, 1, 9.358310000000001, 28.07493, 0.1719944618809179, 0.1719944618809179, 25.89887960894, 233.08991648046, 0.02051007339236019, 0.0205100733923602),
    @("ECs", "Tgfb2", "Acvr1", "FAPs", 3, 1, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3, 1, 26.47935433333333, 79.438063, 0.4866586274141183, 0.4866586274141184, 73.28092465428733, 659.528321888586, 0.0580332881427285, 0.05803328814272851),
    @("ECs", "Tgfb2", "Acvr1", "Inflammatory-Mac", 3, 1, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3, 1, 5.152806666666667, 15.45842, 0.09470237786627494, 0.09470237786627496, 14.26025847702667, 128.34232629324, 0.01129311199457768, 0.01129311199457768),
    @("ECs", "Tgfb2", "Acvr1", "MuSCs", 3, 1, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3, 1, 7.974813333333334, 23.92444, 0.146567460136225, 0.146567460136225, 22.07008855485333, 198.63079699368, 0.01747794278636199, 0.017477942786362),
    @("ECs", "Tgfb2", "Acvr1", "Resolving-Mac", 3, 1, 2.767474, 8.302422, 0.1192484523516842, 0.1192484523516842, 3, 1, 5.445246666666667, 16.33574, 0.1000770727024639, 0.1000770727024639, 15.06957857358667, 135.62620716228, 0.0119340360356558, 0.0119340360356558),
    @("FAPs", "Tgfb2", "Acvr1", "ECs", 3, 1, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3, 1, 9.358310000000001, 28.07493, 0.1719944618809179, 0.1719944618809179, 120.5891113540534, 1085.30200218648, 0.09549801232858651, 0.09549801232858653),
    @("FAPs", "Tgfb2", "Acvr1", "FAPs", 3, 1, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3, 1, 26.47935433333333, 79.438063, 0.4866586274141183, 0.4866586274141184, 341.2070991755742, 3070.863892580168, 0.2702117910795515, 0.2702117910795515),
    @("FAPs", "Tgfb2", "Acvr1", "Inflammatory-Mac", 3, 1, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3, 1, 5.152806666666667, 15.45842, 0.09470237786627494, 0.09470237786627496, 66.39792621879111, 597.58133596912, 0.05258244219096782, 0.05258244219096782),
    @("FAPs", "Tgfb2", "Acvr1", "MuSCs", 3, 1, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3, 1, 7.974813333333334, 23.92444, 0.146567460136225, 0.146567460136225, 102.7616795213156, 924.8551156918401, 0.0813799523658484, 0.08137995236584841),
    @("FAPs", "Tgfb2", "Acvr1", "Resolving-Mac", 3, 1, 12.88577866666667, 38.657336, 0.5552388796954726, 0.5552388796954726, 3, 1, 5.445246666666667, 16.33574, 0.1000770727024639, 0.1000770727024639, 70.16624333207112, 631.49618998864, 0.05556668173051843, 0.05556668173051844),
    @("Inflammatory-Mac", "Tgfb2", "Acvr1", "ECs", 1, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3, 1, 9.358310000000001, 28.07493, 0.1719944618809179, 0.1719944618809179, 0.3203505484833334, 2.88315493635, 0.0002536948840986919, 0.000253694884098692),
    @("Inflammatory-Mac", "Tgfb2", "Acvr1", "FAPs", 1, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3, 1, 26.47935433333333, 79.438063, 0.4866586274141183, 0.4866586274141184, 0.9064324310872223, 8.157891879785, 0.0007178301134075698, 0.00071783011340757),
    @("Inflammatory-Mac", "Tgfb2", "Acvr1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3, 1, 5.152806666666667, 15.45842, 0.09470237786627494, 0.09470237786627496, 0.1763891602111111, 1.5875024419, 0.0001396876882773671, 0.0001396876882773671),
    @("Inflammatory-Mac", "Tgfb2", "Acvr1", "MuSCs", 1, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3, 1, 7.974813333333334, 23.92444, 0.146567460136225, 0.146567460136225, 0.2729911517555556, 2.4569203658, 0.0002161896052074256, 0.0002161896052074257),
    @("Inflammatory-Mac", "Tgfb2", "Acvr1", "Resolving-Mac", 1, 0.3333333333333333, 0.03423166666666667, 0.102695, 0.001475017749550216, 0.001475017749550216, 3, 1, 5.445246666666667, 16.33574, 0.1000770727024639, 0.1000770727024639, 0.1863998688111111, 1.6775988193, 0.0001476154585591617, 0.0001476154585591617),
    @("MuSCs", "Tgfb2", "Acvr1", "ECs", 3, 1, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3, 1, 9.358310000000001, 28.07493, 0.1719944618809179, 0.1719944618809179, 70.37585751326, 633.3827176193399, 0.05573268127587248, 0.0557326812758725),
    @("MuSCs", "Tgfb2", "Acvr1", "FAPs", 3, 1, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3, 1, 26.47935433333333, 79.438063, 0.4866586274141183, 0.4866586274141184, 199.1286105723993, 1792.157495151594, 0.1576957180784308, 0.1576957180784308),
    @("MuSCs", "Tgfb2", "Acvr1", "Inflammatory-Mac", 3, 1, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3, 1, 5.152806666666667, 15.45842, 0.09470237786627494, 0.09470237786627496, 38.74985844310667, 348.74872598796, 0.03068713599245208, 0.03068713599245209),
    @("MuSCs", "Tgfb2", "Acvr1", "MuSCs", 3, 1, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3, 1, 7.974813333333334, 23.92444, 0.146567460136225, 0.146567460136225, 59.97176058941333, 539.74584530472, 0.04749337537880716, 0.04749337537880718),
    @("MuSCs", "Tgfb2", "Acvr1", "Resolving-Mac", 3, 1, 7.520146, 22.560438, 0.324037650203293, 0.3240376502032931, 3, 1, 5.445246666666667, 16.33574, 0.1000770727024639, 0.1000770727024639, 40.94904993934666, 368.54144945412, 0.03242873947773053, 0.03242873947773054)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
